$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Daily data refresh (new case counts) + re-sort by "Casos totales" desc ---
# The source sheet is kept sorted by column B (Casos totales) descending; a handful
# of countries overtook their neighbours after todays update, so those rows swap
# their country name along with picking up fresh (or shifted) statistics.

$ws.Range('A1').Value = 'Datos actualizados a 17 de Junio de 2020 a las 00:51'
$ws.Range('B4').Value = 2206899
$ws.Range('C4').Value = 23949
$ws.Range('D4').Value = 898718
$ws.Range('E4').Value = 1189077
$ws.Range('G4').Value = 821
$ws.Range('H4').Value = 119104
$ws.Range('B11').Value = 237156
$ws.Range('C11').Value = 4164
$ws.Range('D11').Value = 125205
$ws.Range('E11').Value = 104895
$ws.Range('G11').Value = 196
$ws.Range('H11').Value = 7056
$ws.Range('D14').Value = 156232
$ws.Range('E14').Value = 24834
$ws.Range('B24').Value = 76334
$ws.Range('C24').Value = 2801
$ws.Range('D24').Value = 42063
$ws.Range('E24').Value = 32646
$ws.Range('G24').Value = 57
$ws.Range('H24').Value = 1625
$ws.Range('A27').Value = 'Colombia'
$ws.Range('B27').Value = 54931
$ws.Range('C27').Value = 1868
$ws.Range('D27').Value = 20366
$ws.Range('E27').Value = 32764
$ws.Range('G27').Value = 75
$ws.Range('H27').Value = 1801
$ws.Range('A28').Value = 'Suecia'
$ws.Range('B28').Value = 53323
$ws.Range('C28').Value = 940
$ws.Range('D28').Value = 0
$ws.Range('E28').Value = 0
$ws.Range('G28').Value = 48
$ws.Range('H28').Value = 4939
$ws.Range('B52').Value = 17587
$ws.Range('C52').Value = 85
$ws.Range('D52').Value = 15701
$ws.Range('E52').Value = 959
$ws.Range('G52').Value = 2
$ws.Range('H52').Value = 927
$ws.Range('B65').Value = 10111
$ws.Range('C65').Value = 47
$ws.Range('E65').Value = 2422
$ws.Range('B69').Value = 8660
$ws.Range('C69').Value = 13
$ws.Range('E69').Value = 280
$ws.Range('B71').Value = 7740
$ws.Range('C71').Value = 305
$ws.Range('D71').Value = 2820
$ws.Range('E71').Value = 4443
$ws.Range('G71').Value = 9
$ws.Range('H71').Value = 477
$ws.Range('B83').Value = 4299
$ws.Range('C83').Value = 142
$ws.Range('E83').Value = 2341
$ws.Range('A84').Value = 'Gabon'
$ws.Range('B84').Value = 4114
$ws.Range('C84').Value = 81
$ws.Range('D84').Value = 1432
$ws.Range('E84').Value = 2653
$ws.Range('H84').Value = 29
$ws.Range('A85').Value = 'Hungria'
$ws.Range('B85').Value = 4077
$ws.Range('C85').Value = 1
$ws.Range('D85').Value = 2516
$ws.Range('E85').Value = 996
$ws.Range('G85').Value = 2
$ws.Range('H85').Value = 565
$ws.Range('A86').Value = 'Luxemburgo'
$ws.Range('B86').Value = 4075
$ws.Range('C86').Value = 3
$ws.Range('D86').Value = 3933
$ws.Range('E86').Value = 32
$ws.Range('H86').Value = 110
$ws.Range('D127').Value = 816
$ws.Range('E127').Value = 151
$ws.Range('B165').Value = 193
$ws.Range('C165').Value = 6
$ws.Range('D165').Value = 123
$ws.Range('E165').Value = 69
$ws.Range('B193').Value = 29
$ws.Range('C193').Value = 2
$ws.Range('E193').Value = 4
$ws.Range('A206').Value = 'Islas Malvinas'
$ws.Range('A207').Value = 'Groenlandia'
$ws.Range('A210').Value = 'Montserrat'
$ws.Range('D210').Value = 10
$ws.Range('H210').Value = 1
$ws.Range('A211').Value = 'Seychelles'
$ws.Range('D211').Value = 11
$ws.Range('H211').Value = 0
$ws.Range('A213').Value = 'Islas Virgenes Britanicas'
$ws.Range('D213').Value = 7
$ws.Range('H213').Value = 1
$ws.Range('A214').Value = 'Papua Nueva Guinea'
$ws.Range('D214').Value = 8
$ws.Range('H214').Value = 0
